$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Header date
Replace-Text "2024-08-15 Thursday" "2024-08-16 Friday"

# Table cells (order matters: 48÷7= must become 73÷7= before 60÷7= becomes 48÷7=
# to avoid the second replace re-touching the newly written text)
Replace-Text "42÷2=" "62÷5="
Replace-Text "18÷2=" "50÷4="
Replace-Text "21÷3=" "97÷4="
Replace-Text "45÷9=" "55÷3="
Replace-Text "86÷4=" "11÷7="
Replace-Text "57÷3=" "62÷7="
Replace-Text "57÷4=" "53÷5="
Replace-Text "26÷9=" "76÷4="
Replace-Text "90÷3=" "55÷9="
Replace-Text "49÷4=" "12÷7="
Replace-Text "48÷7=" "73÷7="
Replace-Text "60÷7=" "48÷7="
Replace-Text "22÷2=" "56÷2="
Replace-Text "11÷8=" "62÷8="
Replace-Text "44÷5=" "22÷7="
Replace-Text "53÷4=" "19÷2="
Replace-Text "53÷8=" "87÷5="
Replace-Text "43÷2=" "41÷7="
Replace-Text "21÷5=" "42÷3="
Replace-Text "52÷8=" "11÷3="
Replace-Text "53÷3=" "12÷8="
Replace-Text "70÷6=" "76÷6="
Replace-Text "45÷7=" "68÷3="
Replace-Text "87÷3=" "75÷9="
Replace-Text "96÷4=" "62÷9="
